# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.159.76"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.833.96"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.95"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6592"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07411"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07767"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "1.838.64"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6651"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.101"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008550"
$ws.Range("E17").Value = "  +4.52%  "
$ws.Range("D18").Value = "29.170.56"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "2.115.92"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "226.74"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.095"
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.0000"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.59"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.597"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1392"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.94"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.194"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05269"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.864"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7373"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.143"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "1.301.71"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9190"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.034"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08559"
$ws.Range("E43").Value = "  +14.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.60"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "2.023.04"
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5142"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.51"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.750"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05843"
$ws.Range("E51").Value = "  -1.09%  "
